# Revert "Uploading newest EPS-US files"
#
# Restores the older layout of the workbook:
#   About | Calculations (new sheet) | CApULAbIFM
#
# - About: the single "consultation with American Forest Foundation" note
#   is replaced by the full EPA citation (publisher / year / title /
#   hyperlinked URL / page reference).
# - Calculations: new sheet with the low/high/average tons-CO2 figures and
#   the gram conversion that CApULAbIFM!B2 now points at.
# - CApULAbIFM: B2 becomes "=Calculations!A6" instead of the hard-coded
#   "=1.5*10^6" literal.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Free up the shared-string slot used by the old citation before the new
# sheet claims it, so new strings land at the same indices as the target
# workbook.
$about.Range("B3").ClearContents()

# --- Insert the "Calculations" sheet between "About" and "CApULAbIFM" ---
$calc = $wb.Worksheets.Add($null, $about)
$calc.Name = "Calculations"

$calc.Range("A1").Value = "Increased Annual CO2 Sequestration Achievable by Improved Management Practices per Acre"
$calc.Range("A1").Font.Bold = $true

$calc.Range("A2").Value = 2.1
$calc.Range("B2").Value = "tons CO2 / acre / yr"
$calc.Range("C2").Value = "Low Estimate"

$calc.Range("A3").Value = 3.1
$calc.Range("B3").Value = "tons CO2 / acre / yr"
$calc.Range("C3").Value = "High Estimate"

$calc.Range("A4").Formula = "=AVERAGE(A2:A3)"
$calc.Range("B4").Value = "tons CO2 / acre / yr"
$calc.Range("C4").Value = "Average"

$calc.Range("A6").Formula = "=A4*10^6"
$calc.Range("A6").NumberFormat = "0.00E+00"
$calc.Range("B6").Value = "g CO2 / acre / yr"
$calc.Range("C6").Value = "Average, converted to grams CO2"

$calc.Columns("B").ColumnWidth = 18.1666666666667

# --- Update the "About" sheet with the full EPA citation ---
$about.Range("B3").Value = "U.S. EPA"

# Build the "2005" cell with left alignment on A5 first (which mutates the
# pre-existing-but-unused style slot in place), then move it over to B4 -
# this reproduces the original author's formatting history faithfully.
$about.Range("A5").Value = 2005
$about.Range("A5").HorizontalAlignment = -4131  # xlLeft
$about.Range("A5").Cut($about.Range("B4"))
$about.Range("A5").ClearFormats()

$about.Range("B5").Value = "Greenhouse Gas Mitigation Potential in U.S. Forestry and Agriculture"

$about.Range("B6").Value = "http://www.epa.gov/climate/climatechange/Downloads/ccs/ghg_mitigation_forestry_ag_2005.pdf"
$about.Hyperlinks.Add($about.Range("B6"), "http://www.epa.gov/climate/climatechange/Downloads/ccs/ghg_mitigation_forestry_ag_2005.pdf") | Out-Null

$about.Range("B7").Value = "Page 2-3, Table 2-1"

# Column A is bold all the way down alongside the citation rows.
$about.Range("A4").Font.Bold = $true
$about.Range("A5").Font.Bold = $true
$about.Range("A6").Font.Bold = $true
$about.Range("A7").Font.Bold = $true
$about.Range("A9").Font.Bold = $true

# --- Update the "CApULAbIFM" sheet to pull from Calculations ---
$ca = $wb.Worksheets.Item("CApULAbIFM")
$ca.Range("B2").Formula = "=Calculations!A6"

# Leave "About" as the active/selected sheet, as in the target workbook.
$about.Activate()
